$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.638.44"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.671.24"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.05"
$ws.Range("E5").Value = "  -7.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.79"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.496"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.19"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.290.35"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.38"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.715.12"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.694.49"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.52"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.84"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.36"
$ws.Range("E20").Value = "  +5.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.86"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.819.67"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  -3.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.69"
$ws.Range("E28").Value = "  -4.92%  "
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.67"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.98"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.58"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.38"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.672.06"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.27"
$ws.Range("E37").Value = "  -3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "179.78"
$ws.Range("E39").Value = "  +4.57%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0889"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.925"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.75"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.87"
$ws.Range("E46").Value = "  +3.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.84"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000263"
$ws.Range("E49").Value = "  -6.67%  "
$ws.Range("E50").Value = "  -5.46%  "
$ws.Range("E51").Value = "  -6.16%  "
